# Apply "Horarios actualizados Linea 141 - 963" scrape refresh
# Updates header metadata (last-updated timestamp + row count) and the
# changed/new data rows on each of the 3 schedule sheets.
$wb = $excel.ActiveWorkbook

# ----- Sheet: LP1912 -----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,1).Value = "Última actualización: 10:07:51"
$ws.Cells.Item(3,1).Value = "Total filas: 119"

$ws.Cells.Item(16,1).Value = "06:38:54"
$ws.Cells.Item(16,3).Value = "16_SANTA ANA"
$ws.Cells.Item(16,4).Value = 2
$ws.Cells.Item(17,1).Value = "05:44:02"
$ws.Cells.Item(17,3).Value = "17X38_ROMERO"
$ws.Cells.Item(17,4).Value = 56
$ws.Cells.Item(28,1).Value = "07:15:48"
$ws.Cells.Item(28,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(28,4).Value = 6
$ws.Cells.Item(29,1).Value = "06:56:24"
$ws.Cells.Item(29,3).Value = "16_SANTA ANA"
$ws.Cells.Item(29,4).Value = 25
$ws.Cells.Item(41,1).Value = "07:52:32"
$ws.Cells.Item(41,3).Value = "17_ROMERO"
$ws.Cells.Item(41,4).Value = 8
$ws.Cells.Item(42,1).Value = "06:38:54"
$ws.Cells.Item(42,3).Value = "16_SANTA ANA"
$ws.Cells.Item(42,4).Value = 82
$ws.Cells.Item(51,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(52,3).Value = "15_ABASTO"
$ws.Cells.Item(53,3).Value = "10_OLMOS"
$ws.Cells.Item(54,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(64,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(65,3).Value = "215B_EL PATO"
$ws.Cells.Item(78,1).Value = "08:52:33"
$ws.Cells.Item(78,3).Value = "15X38_ABASTO"
$ws.Cells.Item(78,4).Value = 26
$ws.Cells.Item(79,1).Value = "08:30:14"
$ws.Cells.Item(79,3).Value = "14_ABASTO"
$ws.Cells.Item(79,4).Value = 48
$ws.Cells.Item(95,1).Value = "10:07:51"
$ws.Cells.Item(95,4).Value = 6
$ws.Cells.Item(98,1).Value = "10:07:51"
$ws.Cells.Item(98,2).Value = "10:22"
$ws.Cells.Item(98,4).Value = 15
$ws.Cells.Item(99,1).Value = "08:30:14"
$ws.Cells.Item(99,2).Value = "10:24"
$ws.Cells.Item(99,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(99,4).Value = 114
$ws.Cells.Item(100,1).Value = "10:07:51"
$ws.Cells.Item(100,2).Value = "10:25"
$ws.Cells.Item(100,3).Value = "16_SANTA ANA"
$ws.Cells.Item(100,4).Value = 18
$ws.Cells.Item(101,1).Value = "08:40:59"
$ws.Cells.Item(101,2).Value = "10:28"
$ws.Cells.Item(101,3).Value = "15_ABASTO"
$ws.Cells.Item(101,4).Value = 108
$ws.Cells.Item(102,1).Value = "10:07:51"
$ws.Cells.Item(102,3).Value = "14_ABASTO"
$ws.Cells.Item(102,4).Value = 22
$ws.Cells.Item(103,1).Value = "10:07:51"
$ws.Cells.Item(103,2).Value = "10:29"
$ws.Cells.Item(103,3).Value = "15_ABASTO"
$ws.Cells.Item(103,4).Value = 22
$ws.Cells.Item(104,1).Value = "10:07:51"
$ws.Cells.Item(104,2).Value = "10:43"
$ws.Cells.Item(104,4).Value = 36
$ws.Cells.Item(105,1).Value = "08:52:33"
$ws.Cells.Item(105,2).Value = "10:44"
$ws.Cells.Item(105,3).Value = "11X44_ETCHEVERRY"
$ws.Cells.Item(105,4).Value = 112
$ws.Cells.Item(106,1).Value = "10:07:51"
$ws.Cells.Item(106,2).Value = "10:46"
$ws.Cells.Item(106,3).Value = "15_P INDUSTRIAL"
$ws.Cells.Item(106,4).Value = 39
$ws.Cells.Item(107,2).Value = "10:53"
$ws.Cells.Item(107,3).Value = "27_EL RETIRO"
$ws.Cells.Item(107,4).Value = 90
$ws.Cells.Item(108,1).Value = "10:07:51"
$ws.Cells.Item(108,2).Value = "10:56"
$ws.Cells.Item(108,3).Value = "27_EL RETIRO"
$ws.Cells.Item(108,4).Value = 49
$ws.Cells.Item(109,1).Value = "10:07:51"
$ws.Cells.Item(109,2).Value = "10:59"
$ws.Cells.Item(109,3).Value = "10_OLMOS"
$ws.Cells.Item(109,4).Value = 52
$ws.Cells.Item(110,2).Value = "11:01"
$ws.Cells.Item(110,3).Value = "10_OLMOS"
$ws.Cells.Item(110,4).Value = 98
$ws.Cells.Item(111,1).Value = "10:07:51"
$ws.Cells.Item(111,2).Value = "11:01"
$ws.Cells.Item(111,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(111,4).Value = 54
$ws.Cells.Item(112,1).Value = "10:07:51"
$ws.Cells.Item(112,2).Value = "11:03"
$ws.Cells.Item(112,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(112,4).Value = 56
$ws.Cells.Item(112,5).Value = "LP1912"
$ws.Cells.Item(113,1).Value = "10:07:51"
$ws.Cells.Item(113,2).Value = "11:10"
$ws.Cells.Item(113,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(113,4).Value = 63
$ws.Cells.Item(113,5).Value = "LP1912"
$ws.Cells.Item(114,1).Value = "10:07:51"
$ws.Cells.Item(114,2).Value = "11:14"
$ws.Cells.Item(114,3).Value = "14_ABASTO"
$ws.Cells.Item(114,4).Value = 67
$ws.Cells.Item(114,5).Value = "LP1912"
$ws.Cells.Item(115,1).Value = "10:07:51"
$ws.Cells.Item(115,2).Value = "11:15"
$ws.Cells.Item(115,3).Value = "15X38_ABASTO"
$ws.Cells.Item(115,4).Value = 68
$ws.Cells.Item(115,5).Value = "LP1912"
$ws.Cells.Item(116,1).Value = "09:23:52"
$ws.Cells.Item(116,2).Value = "11:19"
$ws.Cells.Item(116,3).Value = "10_OLMOS"
$ws.Cells.Item(116,4).Value = 116
$ws.Cells.Item(116,5).Value = "LP1912"
$ws.Cells.Item(117,1).Value = "10:07:51"
$ws.Cells.Item(117,2).Value = "11:21"
$ws.Cells.Item(117,3).Value = "16_SANTA ANA"
$ws.Cells.Item(117,4).Value = 74
$ws.Cells.Item(117,5).Value = "LP1912"
$ws.Cells.Item(118,1).Value = "10:07:51"
$ws.Cells.Item(118,2).Value = "11:29"
$ws.Cells.Item(118,3).Value = "10_OLMOS"
$ws.Cells.Item(118,4).Value = 82
$ws.Cells.Item(118,5).Value = "LP1912"
$ws.Cells.Item(119,1).Value = "10:07:51"
$ws.Cells.Item(119,2).Value = "11:30"
$ws.Cells.Item(119,3).Value = "215C_EL PATO"
$ws.Cells.Item(119,4).Value = 83
$ws.Cells.Item(119,5).Value = "LP1912"
$ws.Cells.Item(120,1).Value = "10:07:51"
$ws.Cells.Item(120,2).Value = "11:41"
$ws.Cells.Item(120,3).Value = "215B_EL PATO"
$ws.Cells.Item(120,4).Value = 94
$ws.Cells.Item(120,5).Value = "LP1912"
$ws.Cells.Item(121,1).Value = "10:07:51"
$ws.Cells.Item(121,2).Value = "11:45"
$ws.Cells.Item(121,3).Value = "15X38_ABASTO"
$ws.Cells.Item(121,4).Value = 98
$ws.Cells.Item(121,5).Value = "LP1912"
$ws.Cells.Item(122,1).Value = "10:07:51"
$ws.Cells.Item(122,2).Value = "11:52"
$ws.Cells.Item(122,3).Value = "225_GOMEZ"
$ws.Cells.Item(122,4).Value = 105
$ws.Cells.Item(122,5).Value = "LP1912"
$ws.Cells.Item(123,1).Value = "10:07:51"
$ws.Cells.Item(123,2).Value = "11:58"
$ws.Cells.Item(123,3).Value = "17_ROMERO"
$ws.Cells.Item(123,4).Value = 111
$ws.Cells.Item(123,5).Value = "LP1912"
$ws.Cells.Item(124,1).Value = "10:07:51"
$ws.Cells.Item(124,2).Value = "12:05"
$ws.Cells.Item(124,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(124,4).Value = 118
$ws.Cells.Item(124,5).Value = "LP1912"

# ----- Sheet: LP1912-215 -----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,1).Value = "Última actualización: 10:07:51"
$ws.Cells.Item(3,1).Value = "Total filas: 17"

$ws.Cells.Item(21,1).Value = "10:07:51"
$ws.Cells.Item(21,2).Value = "11:30"
$ws.Cells.Item(21,3).Value = "215C_EL PATO"
$ws.Cells.Item(21,4).Value = 83
$ws.Cells.Item(21,5).Value = "LP1912"
$ws.Cells.Item(22,1).Value = "10:07:51"
$ws.Cells.Item(22,2).Value = "11:41"
$ws.Cells.Item(22,3).Value = "215B_EL PATO"
$ws.Cells.Item(22,4).Value = 94
$ws.Cells.Item(22,5).Value = "LP1912"

# ----- Sheet: 6203-6173 -----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,1).Value = "Última actualización: 10:07:51"
$ws.Cells.Item(3,1).Value = "Total filas: 17"

$ws.Cells.Item(16,1).Value = "10:07:51"
$ws.Cells.Item(16,4).Value = 5
$ws.Cells.Item(18,1).Value = "10:07:51"
$ws.Cells.Item(18,4).Value = 22
$ws.Cells.Item(19,1).Value = "08:52:33"
$ws.Cells.Item(19,3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(19,4).Value = 98
$ws.Cells.Item(20,1).Value = "10:07:51"
$ws.Cells.Item(20,3).Value = "215A_LA PLATA"
$ws.Cells.Item(20,4).Value = 23
$ws.Cells.Item(22,1).Value = "10:07:51"
$ws.Cells.Item(22,2).Value = "11:25"
$ws.Cells.Item(22,3).Value = "215C_LA PLATA"
$ws.Cells.Item(22,4).Value = 78
$ws.Cells.Item(22,5).Value = "L6203"

